$d = $word.ActiveDocument

$d.Content.Find.Execute("54×69=3726", $true, $false, $false, $false, $false, $true, 1, $false, "31×38=1178", 2) | Out-Null
$d.Content.Find.Execute("51×78=3978", $true, $false, $false, $false, $false, $true, 1, $false, "74×84=6216", 2) | Out-Null
$d.Content.Find.Execute("77×75=5775", $true, $false, $false, $false, $false, $true, 1, $false, "99×43=4257", 2) | Out-Null
$d.Content.Find.Execute("45×12=540", $true, $false, $false, $false, $false, $true, 1, $false, "13×30=390", 2) | Out-Null
$d.Content.Find.Execute("70×34=2380", $true, $false, $false, $false, $false, $true, 1, $false, "22×69=1518", 2) | Out-Null
$d.Content.Find.Execute("82×77=6314", $true, $false, $false, $false, $false, $true, 1, $false, "99×69=6831", 2) | Out-Null
$d.Content.Find.Execute("17×77=1309", $true, $false, $false, $false, $false, $true, 1, $false, "38×47=1786", 2) | Out-Null
$d.Content.Find.Execute("89×73=6497", $true, $false, $false, $false, $false, $true, 1, $false, "77×23=1771", 2) | Out-Null
$d.Content.Find.Execute("54×31=1674", $true, $false, $false, $false, $false, $true, 1, $false, "78×55=4290", 2) | Out-Null
$d.Content.Find.Execute("94×91=8554", $true, $false, $false, $false, $false, $true, 1, $false, "24×77=1848", 2) | Out-Null
$d.Content.Find.Execute("99×74=7326", $true, $false, $false, $false, $false, $true, 1, $false, "16×56=896", 2) | Out-Null
$d.Content.Find.Execute("94×39=3666", $true, $false, $false, $false, $false, $true, 1, $false, "94×53=4982", 2) | Out-Null
$d.Content.Find.Execute("73×65=4745", $true, $false, $false, $false, $false, $true, 1, $false, "52×47=2444", 2) | Out-Null
$d.Content.Find.Execute("34×60=2040", $true, $false, $false, $false, $false, $true, 1, $false, "86×79=6794", 2) | Out-Null
$d.Content.Find.Execute("17×67=1139", $true, $false, $false, $false, $false, $true, 1, $false, "12×52=624", 2) | Out-Null
$d.Content.Find.Execute("28×24=672", $true, $false, $false, $false, $false, $true, 1, $false, "37×56=2072", 2) | Out-Null
$d.Content.Find.Execute("21×56=1176", $true, $false, $false, $false, $false, $true, 1, $false, "27×55=1485", 2) | Out-Null
$d.Content.Find.Execute("66×60=3960", $true, $false, $false, $false, $false, $true, 1, $false, "64×96=6144", 2) | Out-Null
$d.Content.Find.Execute("11×67=737", $true, $false, $false, $false, $false, $true, 1, $false, "95×86=8170", 2) | Out-Null
$d.Content.Find.Execute("99×75=7425", $true, $false, $false, $false, $false, $true, 1, $false, "59×27=1593", 2) | Out-Null
$d.Content.Find.Execute("67×45=3015", $true, $false, $false, $false, $false, $true, 1, $false, "28×63=1764", 2) | Out-Null
$d.Content.Find.Execute("54×53=2862", $true, $false, $false, $false, $false, $true, 1, $false, "68×18=1224", 2) | Out-Null
$d.Content.Find.Execute("52×98=5096", $true, $false, $false, $false, $false, $true, 1, $false, "32×36=1152", 2) | Out-Null
$d.Content.Find.Execute("90×30=2700", $true, $false, $false, $false, $false, $true, 1, $false, "33×32=1056", 2) | Out-Null
$d.Content.Find.Execute("84×94=7896", $true, $false, $false, $false, $false, $true, 1, $false, "75×29=2175", 2) | Out-Null
